$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row updates to the cryptocurrency price table (cols B-E).
# D-column numeric-looking text values are prefixed with a leading
# apostrophe so Excel keeps them as literal text (matching the source
# data, which stores prices/caps as text, not numbers).

$ws.Range("D2").Value = "58.730.58"
$ws.Range("E2").Value = "  +2.14%  "

$ws.Range("D3").Value = "3.086.50"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'520.46"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").Value = "'143.71"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'0.438"
$ws.Range("E8").Value = "  +0.53%  "

$ws.Range("D9").Value = "'7.34"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").Value = "'0.109"
$ws.Range("E10").Value = "  +0.36%  "

$ws.Range("D11").Value = "'0.383"
$ws.Range("E11").Value = "  +2.28%  "

$ws.Range("D12").Value = "3.619.31"
$ws.Range("E12").Value = "  +0.68%  "

$ws.Range("E13").Value = "  +0.90%  "

$ws.Range("D14").Value = "'26.69"
$ws.Range("E14").Value = "  +3.69%  "

$ws.Range("D15").Value = "'0.0000166"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").Value = "58.718.09"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("D17").Value = "3.094.78"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("D18").Value = "'6.13"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").Value = "'12.87"
$ws.Range("E19").Value = "  -1.30%  "

$ws.Range("D20").Value = "'8.09"
$ws.Range("E20").Value = "  -1.40%  "

$ws.Range("D21").Value = "'343.95"
$ws.Range("E21").Value = "  +1.89%  "

$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").Value = "'0.505"
$ws.Range("E23").Value = "  +0.85%  "

$ws.Range("D24").Value = "'65.65"
$ws.Range("E24").Value = "  +0.32%  "

$ws.Range("D25").Value = "'0.171"
$ws.Range("E25").Value = "  -0.72%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").Value = "0.0₃0920"
$ws.Range("E27").Value = "  -1.35%  "

$ws.Range("D28").Value = "'6.62"
$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("D29").Value = "'7.23"
$ws.Range("E29").Value = "  +1.98%  "

$ws.Range("D30").Value = "'1.84"
$ws.Range("E30").Value = "  +1.74%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'20.98"
$ws.Range("E31").Value = "  +0.61%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.20"
$ws.Range("E32").Value = "  +1.71%  "

$ws.Range("D33").Value = "'154.59"
$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("D34").Value = "'4.61"
$ws.Range("E34").Value = "  +1.76%  "

$ws.Range("D35").Value = "'6.12"
$ws.Range("E35").Value = "  +3.40%  "

$ws.Range("D36").Value = "'26.73"
$ws.Range("E36").Value = "  +0.54%  "

$ws.Range("D37").Value = "'1.29"
$ws.Range("E37").Value = "  +4.06%  "

$ws.Range("D38").Value = "'0.0685"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("B39").Value = "RenzoRestakedETH"
$ws.Range("C39").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D39").Value = "3.129.43"
$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'3.91"
$ws.Range("E40").Value = "  +0.97%  "

$ws.Range("D41").Value = "'36.76"
$ws.Range("E41").Value = "  -0.49%  "

$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("D43").Value = "'0.664"
$ws.Range("E43").Value = "  -0.90%  "

$ws.Range("D44").Value = "'1.45"
$ws.Range("E44").Value = "  +3.96%  "

$ws.Range("D45").Value = "2.285.68"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("D46").Value = "'0.0254"
$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").Value = "'20.72"
$ws.Range("E47").Value = "  +1.83%  "

$ws.Range("D48").Value = "'0.962"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").Value = "'5.98"
$ws.Range("E49").Value = "  +1.79%  "

$ws.Range("D50").Value = "'0.751"
$ws.Range("E50").Value = "  +8.71%  "

$ws.Range("D51").Value = "'262.69"
$ws.Range("E51").Value = "  +11.50%  "
